# "I was added to db" - add a new person (Nurmukhanbet Rakhimbayev) to the
# DB worksheet, rename the "Name" header to "First Name", align the font
# used by the existing name/last name/username cells with the new person's
# font, and extend the 0.00 number formatting down through the rest of the
# Debt/Check/Approve columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header rename: "Name" -> "First Name"
$ws.Range("C1").Value = "First Name"

# 2. New row of data for row 4 (person #3)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 287100650
$ws.Range("C4").Value = "Nurmukhanbet"
$ws.Range("D4").Value = "Rakhimbayev"
$ws.Range("E4").Value = "Nurmukhanbet Rakhimbayev "
$ws.Range("F4").Value = 350
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0

# Match the style already used on the existing data rows (2-3) so the new
# row renders identically: bordered/centered cells, ID column using the
# "name" font, numeric Debt/Check/Approve columns using the 0.00 format.
$ws.Range("A4").Style = $ws.Range("A3").Style
$ws.Range("B4").Style = $ws.Range("B3").Style
$ws.Range("C4:E4").Style = $ws.Range("C3:E3").Style
$ws.Range("F4:H4").Style = $ws.Range("F3:H3").Style

# 3. The new row's Name/Last Name/Username font is Ubuntu; apply that same
# font to the existing rows 2-3 so the whole column is consistent (this is
# what collapses the old dedicated "JetBrains Mono" font into the font
# already used elsewhere in the sheet).
$ws.Range("C4:E4").Font.Name = "Ubuntu"
$ws.Range("C2:E3").Font.Name = "Ubuntu"

# 4. Extend the 0.00 number format used by the populated Debt/Check/Approve
# cells (rows 2-4) down across the remaining empty rows (5-23).
$ws.Range("F5:H23").NumberFormat = "0.00"

# 5. Minor cosmetic leftovers from the editing session.
$ws.Columns("D").ColumnWidth = 23.07
$ws.Range("D12").Select
